$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New TCP (DropTail) values - column B, rows 3-21
$ws.Range("B3").Value = 73.680941176499999
$ws.Range("B4").Value = 64.992000000000004
$ws.Range("B5").Value = 67.248000000000005
$ws.Range("B6").Value = 65.361999999999995
$ws.Range("B7").Value = 65.37
$ws.Range("B8").Value = 65.335999999999999
$ws.Range("B9").Value = 65.325999999999993
$ws.Range("B10").Value = 65.31
$ws.Range("B11").Value = 65.293999999999997
$ws.Range("B12").Value = 65.286000000000001
$ws.Range("B13").Value = 65.262
$ws.Range("B14").Value = 65.254000000000005
$ws.Range("B15").Value = 65.248000000000005
$ws.Range("B16").Value = 65.249610062900004
$ws.Range("B17").Value = 65.249610062900004
$ws.Range("B18").Value = 65.241559748399993
$ws.Range("B19").Value = 65.249610062900004
$ws.Range("B20").Value = 65.249610062900004
$ws.Range("B21").Value = 65.248000000000005

# New CBR (DropTail) values - column C, rows 5-21
$ws.Range("C5").Value = 44.356094276100002
$ws.Range("C6").Value = 43.651884983999999
$ws.Range("C7").Value = 43.6248205128
$ws.Range("C8").Value = 43.624281150199998
$ws.Range("C9").Value = 43.624615384599998
$ws.Range("C10").Value = 43.623258785899999
$ws.Range("C11").Value = 43.622564102600002
$ws.Range("C12").Value = 43.624281150199998
$ws.Range("C13").Value = 43.6235897436
$ws.Range("C14").Value = 43.618146964899999
$ws.Range("C15").Value = 43.626666666699997
$ws.Range("C16").Value = 43.625303514400002
$ws.Range("C17").Value = 43.626666666699997
$ws.Range("C18").Value = 43.625303514400002
$ws.Range("C19").Value = 43.626666666699997
$ws.Range("C20").Value = 43.625303514400002
$ws.Range("C21").Value = 43.2

# New TCP (RED) values - column D, rows 3-21
$ws.Range("D3").Value = 73.680941176499999
$ws.Range("D4").Value = 64.992000000000004
$ws.Range("D5").Value = 67.444313725499995
$ws.Range("D6").Value = 66.4769655172
$ws.Range("D7").Value = 67.054159291999994
$ws.Range("D8").Value = 65.766233766200003
$ws.Range("D9").Value = 65.356276729599998
$ws.Range("D10").Value = 65.342379746800006
$ws.Range("D11").Value = 65.316025157200002
$ws.Range("D12").Value = 65.316025157200002
$ws.Range("D13").Value = 65.281999999999996
$ws.Range("D14").Value = 65.283823899400005
$ws.Range("D15").Value = 65.259672956000003
$ws.Range("D16").Value = 65.251622641500006
$ws.Range("D17").Value = 65.249610062900004
$ws.Range("D18").Value = 65.241559748399993
$ws.Range("D19").Value = 65.249610062900004
$ws.Range("D20").Value = 65.249610062900004
$ws.Range("D21").Value = 65.293176470600002

# New CBR (RED) values - column E, rows 5-21
$ws.Range("E5").Value = 43.496608108099998
$ws.Range("E6").Value = 43.485035143799998
$ws.Range("E7").Value = 44.051897435900003
$ws.Range("E8").Value = 43.751872204500003
$ws.Range("E9").Value = 43.610256410300003
$ws.Range("E10").Value = 43.6171246006
$ws.Range("E11").Value = 43.6108717949
$ws.Range("E12").Value = 43.614670926499997
$ws.Range("E13").Value = 43.630358974400004
$ws.Range("E14").Value = 43.6136485623
$ws.Range("E15").Value = 43.623179487199998
$ws.Range("E16").Value = 43.621827476
$ws.Range("E17").Value = 43.619076923100003
$ws.Range("E18").Value = 43.630006389800002
$ws.Range("E19").Value = 43.619076923100003
$ws.Range("E20").Value = 43.625230769200002
$ws.Range("E21").Value = 43.231897763600003

# Update selection to match target
$ws.Range("D3:D21").Select()
